# Weekly price update: insert a new price record at the top of the
# Berenjena / Macroferia Regional de Talca data block (row 171), pushing
# the existing historical rows (old 171-187) down by one row (172-188).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 171; Excel shifts rows 171:187 down to 172:188
# and extends the used range accordingly.
$ws.Rows.Item(171).Insert()

# Populate the new row 171 with the new weekly record.
$ws.Cells.Item(171, 1).Value = 5
$ws.Cells.Item(171, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(171, 3).Value = "Maule"
$ws.Cells.Item(171, 4).Value = 45132
$ws.Cells.Item(171, 5).Value = 7
$ws.Cells.Item(171, 6).Value = 100112001
$ws.Cells.Item(171, 7).Value = "Berenjena"
$ws.Cells.Item(171, 8).Value = "Sin especificar"
$ws.Cells.Item(171, 9).Value = "Primera"
$ws.Cells.Item(171, 10).Value = 200
$ws.Cells.Item(171, 11).Value = 7000
$ws.Cells.Item(171, 12).Value = 7000
$ws.Cells.Item(171, 13).Value = 7000
$ws.Cells.Item(171, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(171, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(171, 16).Value = 140
$ws.Cells.Item(171, 17).Value = 50
$ws.Cells.Item(171, 18).Value = "Hortaliza"
